$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.096.07"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "1.892.50"
$ws.Range("E3").Value = "  +1.75%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'306.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("D6").Value = "'0.9988"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "'0.5147"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.04%  "

$ws.Range("D8").Value = "'0.3755"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.32%  "

$ws.Range("D9").Value = "'0.07207"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("D10").Value = "'21.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.46%  "

$ws.Range("D11").Value = "'0.9052"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("D12").Value = "'0.07641"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.02%  "

$ws.Range("D13").Value = "1.882.66"
$ws.Range("E13").Value = "  +1.09%  "

$ws.Range("D14").Value = "'95.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.06%  "

$ws.Range("D15").Value = "'5.266"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").Value = "'0.9989"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").Value = "'0.000008487"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").Value = "'14.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.20%  "

$ws.Range("D19").Value = "'0.9984"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").Value = "27.130.46"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("D22").Value = "2.137.54"
$ws.Range("E22").Value = "  +1.55%  "

$ws.Range("D23").Value = "'10.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.04%  "

$ws.Range("D24").Value = "'6.410"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("D25").Value = "'2.295"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.23%  "

$ws.Range("D26").Value = "'145.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.29%  "

$ws.Range("D27").Value = "'1.768"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").Value = "'18.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("D29").Value = "'114.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.26%  "

$ws.Range("D30").Value = "'4.951"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.73%  "

$ws.Range("D31").Value = "'4.827"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.06%  "

$ws.Range("D32").Value = "'0.09186"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.82%  "

$ws.Range("D33").Value = "'0.05084"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "

$ws.Range("D34").Value = "'1.238"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.63%  "

$ws.Range("E35").Value = "  +4.43%  "

$ws.Range("D36").Value = "'2.979"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("D38").Value = "'2.619"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.74%  "

$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").Value = "'0.5598"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'1.077"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").Value = "'9.101"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.21%  "

$ws.Range("D43").Value = "'6.669"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.68%  "

$ws.Range("D44").Value = "'117.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.28%  "

$ws.Range("E45").Value = "  +2.80%  "

$ws.Range("D46").Value = "'0.4815"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.11%  "

$ws.Range("E47").Value = "  +1.27%  "

$ws.Range("D48").Value = "'0.9984"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").Value = "'1.598"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.01%  "

$ws.Range("D50").Value = "'37.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.75%  "

$ws.Range("D51").Value = "'64.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.55%  "
